# Fix the "Reference" column in Sheet1: each entry in column A (rows 2-45)
# has an erroneous trailing "16" appended to the verse number
# (e.g. "Jonah 1:116" should read "Jonah 1:1"). Strip that stray suffix so
# the references are human readable again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $value = $cell.Value2
    if ($value -ne $null -and $value.ToString().EndsWith("16")) {
        $cell.Value2 = $value.ToString().Substring(0, $value.ToString().Length - 2)
    }
}
